$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, reusing the same formatting as the other
# header cells (E1) via copy/paste-special (keeps the existing style
# index instead of minting a new one).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Add data cells F2 and F3 with time_taken values
$ws.Range("F2").Value = "2021-10-05 13:39:08.797879"
$ws.Range("F3").Value = "2021-10-05 13:39:08.797891"
